$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit swaps the data of row 4 <-> row 5, and row 6 <-> row 7
# (the Id/TaxonId/Artnamn/.../Ost/Nord values move between the two rows
# in each pair, while every other column on those rows already holds an
# identical value across the pair and therefore needs no change).

function Set-EmptyTextCell($addr) {
    # Produces an empty (zero-length) TEXT cell - as opposed to a
    # genuinely blank/absent cell - by using Excel's text-prefix trick,
    # then stripping the quote-prefix style it leaves behind so no
    # formatting diff is introduced.
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).ClearFormats()
}

function Clear-Cell($addr) {
    # Fully removes a cell's content (value and "is present" marker).
    $ws.Range($addr).ClearContents()
}

# --- Row 4 <-> Row 5 -------------------------------------------------
$a4 = $ws.Range("A4").Value(); $a5 = $ws.Range("A5").Value()
$b4 = $ws.Range("B4").Value(); $b5 = $ws.Range("B5").Value()
$e4 = $ws.Range("E4").Value(); $e5 = $ws.Range("E5").Value()
$f4 = $ws.Range("F4").Value(); $f5 = $ws.Range("F5").Value()
$g4 = $ws.Range("G4").Value(); $g5 = $ws.Range("G5").Value()
$h4 = $ws.Range("H4").Value(); $h5 = $ws.Range("H5").Value()
$q4 = $ws.Range("Q4").Value(); $q5 = $ws.Range("Q5").Value()
$r4 = $ws.Range("R4").Value(); $r5 = $ws.Range("R5").Value()

$ws.Range("A4").Value = $a5
$ws.Range("B4").Value = $b5
$ws.Range("E4").Value = $e5
$ws.Range("F4").Value = $f5
$ws.Range("G4").Value = $g5
$ws.Range("H4").Value = $h5
$ws.Range("Q4").Value = $q5
$ws.Range("R4").Value = $r5

$ws.Range("A5").Value = $a4
$ws.Range("B5").Value = $b4
$ws.Range("E5").Value = $e4
$ws.Range("F5").Value = $f4
$ws.Range("G5").Value = $g4
$ws.Range("H5").Value = $h4
$ws.Range("Q5").Value = $q4
$ws.Range("R5").Value = $r4

# L4 goes from "no cell" to an empty text cell; L5 the opposite.
Set-EmptyTextCell "L4"
Clear-Cell "L5"

# --- Row 6 <-> Row 7 -------------------------------------------------
$a6 = $ws.Range("A6").Value(); $a7 = $ws.Range("A7").Value()
$b6 = $ws.Range("B6").Value(); $b7 = $ws.Range("B7").Value()
$e6 = $ws.Range("E6").Value(); $e7 = $ws.Range("E7").Value()
$f6 = $ws.Range("F6").Value(); $f7 = $ws.Range("F7").Value()
$g6 = $ws.Range("G6").Value(); $g7 = $ws.Range("G7").Value()
$h6 = $ws.Range("H6").Value(); $h7 = $ws.Range("H7").Value()
$q6 = $ws.Range("Q6").Value(); $q7 = $ws.Range("Q7").Value()
$r6 = $ws.Range("R6").Value(); $r7 = $ws.Range("R7").Value()
$m6 = $ws.Range("M6").Value()

$ws.Range("A6").Value = $a7
$ws.Range("B6").Value = $b7
$ws.Range("E6").Value = $e7
$ws.Range("F6").Value = $f7
$ws.Range("G6").Value = $g7
$ws.Range("H6").Value = $h7
$ws.Range("Q6").Value = $q7
$ws.Range("R6").Value = $r7

$ws.Range("A7").Value = $a6
$ws.Range("B7").Value = $b6
$ws.Range("E7").Value = $e6
$ws.Range("F7").Value = $f6
$ws.Range("G7").Value = $g6
$ws.Range("H7").Value = $h6
$ws.Range("Q7").Value = $q6
$ws.Range("R7").Value = $r6

# J6 goes from "no cell" to an empty text cell; J7 the opposite.
Set-EmptyTextCell "J6"
Clear-Cell "J7"

# M6 ("färska spår") moves to M7; M6 becomes "no cell".
$ws.Range("M7").Value = $m6
Clear-Cell "M6"

# AF6 goes from "no cell" to an empty text cell; AF7 the opposite.
Set-EmptyTextCell "AF6"
Clear-Cell "AF7"

Write-Output "done"
